$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 43542.059783935547
$ws.Range("C2").Value = 82756.7265625
$ws.Range("D2").Value = 52.614524841308594

$ws.Range("B3").Value = 43457.301040649414
$ws.Range("C3").Value = 82756.7265625
$ws.Range("D3").Value = 52.512107849121094

$ws.Range("B4").Value = 70495.849090576172
$ws.Range("C4").Value = 82756.7265625
$ws.Range("D4").Value = 85.184432983398438

$ws.Range("B5").Value = 82756.722778320313
$ws.Range("C5").Value = 82756.7265625
$ws.Range("D5").Value = 99.999992370605469

$ws.Range("B6").Value = 85657.906433105469
$ws.Range("C6").Value = 82756.7265625
$ws.Range("D6").Value = 103.50567626953125
